$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Primera) - date moves from 44223 to 44559, values updated
$ws.Range("D2").Value = 44559
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6500
$ws.Range("S2").Value = 3250

# Row 3 (Segunda) - date moves from 44223 to 44559, values updated
$ws.Range("D3").Value = 44559
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("S3").Value = 2500

# Row 4 - date moves from 44559 to 44574
$ws.Range("D4").Value = 44574

# Row 5 - date moves from 44559 to 44574
$ws.Range("D5").Value = 44574

# Row 6 (Primera) - date moves from 44574 to 44223, values updated
$ws.Range("D6").Value = 44223
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 3500
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 3750
$ws.Range("S6").Value = 1875

# Row 7 (Segunda) - date moves from 44574 to 44223, values updated
$ws.Range("D7").Value = 44223
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("S7").Value = 1500
